# Standardize some naming and formatting
# Replace the text genotype labels "36625-8", "36625-10", "36625-14" in
# column A with their plain numeric identifiers (8, 10, 14), and move the
# active selection to A44 to match where the edits were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14-25 were labelled "36625-8" -> now just the number 8
$ws.Range("A14:A25").Value = 8

# Rows 26-37 were labelled "36625-10" -> now just the number 10
$ws.Range("A26:A37").Value = 10

# Rows 38-49 were labelled "36625-14" -> now just the number 14
$ws.Range("A38:A49").Value = 14

# Update the current selection to reflect where editing last occurred
$ws.Range("A44").Select()
